$d = $word.ActiveDocument

# Namespace declarations reused for the InsertXML payloads below.
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Set-ParagraphXml($paragraphIndex, $expectedSnippet, $pXml) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $r = $p.Range
    # Paragraph.Range.Text includes the trailing paragraph-mark (CR), so
    # match with a trailing wildcard rather than an exact comparison.
    if ($r.Text -notlike ($expectedSnippet + "*")) {
        throw "Paragraph $paragraphIndex did not match expected text: [$($r.Text)]"
    }
    $xmlFrag = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document $ns><w:body>$pXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $r.InsertXML($xmlFrag)
}

# ---------------------------------------------------------------------------
# 1) "...before the client even requests for it, thus reducing load times."
#    -> "...before the client even requests it, thus reducing load times."
# ---------------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "the client even requests for it, thus reducing load times.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the client even requests it, thus reducing load times.", 2)
Write-Host "Change 1 (server push wording) applied: $found1"

# ---------------------------------------------------------------------------
# 2) Heading "HTT" + "P Status Codes" (split across two runs, wrapped in a
#    proofErr spellStart/spellEnd pair) -> single run "HTTP Status Codes"
#    with the proofErr markers removed entirely.
# ---------------------------------------------------------------------------
$pHttpStatusCodes = '<w:p w14:paraId="371D5C98" w14:textId="6F3A464B" w:rsidR="0043732F" w:rsidRDefault="0043732F" w:rsidP="00E867C0"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr><w:t>HTTP Status Codes</w:t></w:r></w:p>'

Set-ParagraphXml 53 "HTTP Status Codes" $pHttpStatusCodes
Write-Host "Change 2 (HTTP Status Codes heading) applied"


# ---------------------------------------------------------------------------
# 3) "201 Created: ... being created (e.g. after a POST)."
#    -> split into three runs, inserting a comma: "(e.g., after a POST)."
# ---------------------------------------------------------------------------
$p201 = '<w:p w14:paraId="65F6DD2A" w14:textId="40FC1230" w:rsidR="003942D8" w:rsidRDefault="003942D8" w:rsidP="003942D8"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr><w:t>201 Created: Successful request that results in a new resource being created (e.g.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> after a POST).</w:t></w:r></w:p>'

Set-ParagraphXml 59 "201 Created: Successful request that results in a new resource being created (e.g. after a POST)." $p201
Write-Host "Change 3 (201 Created comma) applied"

# ---------------------------------------------------------------------------
# 4) "401 Unauthorized: Authentication is required (or has failed)."
#    -> "401 Unauthorized: Authorization is required (or has failed)."
#    split into three runs.
# ---------------------------------------------------------------------------
$p401 = '<w:p w14:paraId="216D0DDF" w14:textId="58034319" w:rsidR="00D9767D" w:rsidRDefault="00D9767D" w:rsidP="00D9767D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">401 Unauthorized: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr><w:t>Authorization</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> is required (or has failed).</w:t></w:r></w:p>'

Set-ParagraphXml 67 "401 Unauthorized: Authentication is required (or has failed)." $p401
Write-Host "Change 4 (401 Unauthorized -> Authorization) applied"

Write-Host "All changes applied."
